# Add new "Joueurs" rows to the BDD workbook (persistance cote serveur).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Joueurs")

$rows = @(
    @(3,  "jacques@hotmail.fr",  "leclerc", "jacques", "jac"),
    @(4,  "jonhattan@gmail.com", "Dupont",  "Jonathan", "jojo"),
    @(5,  "jonhattan@gmail.com", "J",       "jjs",      "jojo"),
    @(6,  "a",                   "a",       "a",        "a"),
    @(7,  "jonhattan@gmail.com", "fb",      "zb",       "uy"),
    @(8,  "jonhattan@gmail.com", "njk",     "mgu",      "kjnj"),
    @(9,  "jonhattan@gmail.com", "vuy",     "jkb",      "kjjb"),
    @(10, "jonhattan@gmail.com", "sdf",     "eth",      "rjjy"),
    @(11, "jacques@hotmail.fr",  "ver",     "umy",      "rhser"),
    @(12, "jacques@hotmail.fr",  "obtrs",   "vbiau",    "cevwg")
)

$r = 4
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}
